# Arbeitszeiten.xlsx update:
# - Add a new time-tracking entry in row 47 (Daniel, "Bulk Upload Excel", 2.5 hours,
#   dated 2018-10-29 / serial 43402). This introduces a new shared string
#   "Bulk Upload Excel" and bumps the H3 SUMIF total for "Daniel" from 69.5 to 72.
# - Leave the final selection on cell E45, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date cell A47: copy the existing date formatting from A46 (numFmt "m/d/yyyy")
# so the new cell reuses the same style index instead of creating a new one,
# then set its value to the date serial for 2018-10-29.
[void]$ws.Range("A46").Copy()
[void]$ws.Range("A47").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A47").Value = 43402

# Name, description and duration for the new row.
$ws.Range("B47").Value = "Daniel"
$ws.Range("C47").Value = "Bulk Upload Excel"
$ws.Range("D47").Value = 2.5

# Match the saved selection state in the sheet view.
[void]$ws.Range("E45").Select()
